# Add new column S (year 2022) data and update select precision values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: year header
$ws.Range("S4").Value = 2022

# Row 5 values (updated P/Q/R + new S)
$ws.Range("P5").Value = 23.111083656771282
$ws.Range("Q5").Value = 24.08077930418019
$ws.Range("R5").Value = 19.336931533747723
$ws.Range("S5").Value = 13.600365850576139

# Row 6 values (updated P/Q/R + new S)
$ws.Range("P6").Value = 14.322631450320875
$ws.Range("Q6").Value = 13.073459110725862
$ws.Range("R6").Value = 10.464141365743002
$ws.Range("S6").Value = 9.2742414863791556

# Row 7 values (updated P + new S)
$ws.Range("P7").Value = 23.612622725489956
$ws.Range("S7").Value = 17.303523954725925

# Row 8: new S value
$ws.Range("S8").Value = 205.5

# Copy styles from column R to column S for rows 3-8 (mirrors the "2021" column)
$ws.Range("R3:R8").Copy()
$ws.Range("S3:S8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Update selection
$ws.Range("Q15").Select()
